{"js": "// The note about the old LAN IP (\"192.168.0.14\") is replaced with an\n// updated reminder to check the winnpysoft.com domain / iTrac app.pos\n// media uploads. The new sentence is split across several runs with\n// <w:proofErr> spell-check markers, mirroring what Word's editor would\n// generate while typing the flagged words \"iTrac\" and \"app.pos\".\n\nconst body = context.document.body;\n\n// Locate the paragraph that still references the old IP address.\nconst results = body.search(\"192.168.0.14\", { matchCase: false, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the '192.168.0.14' paragraph to update.\");\n}\n\nconst target = results.items[0].paragraphs.getFirst();\n\n// Grab the paragraph's current OOXML so we can keep its own opening\n// <w:p> attributes (paraId/rsid/...) and <w:pPr> (list style/numbering)\n// exactly as-is; only the runs inside the paragraph change.\nconst ooxmlResult = target.getOoxml();\nawait context.sync();\nconst existingOoxml = ooxmlResult.value;\n\nconst pOpenMatch = existingOoxml.match(/<w:p(\\s[^>]*)?>/);\nconst pOpenTag = pOpenMatch ? pOpenMatch[0] : \"<w:p>\";\nconst pPrMatch = existingOoxml.match(/<w:pPr>[\\s\\S]*?<\\/w:pPr>/);\nconst pPrXml = pPrMatch ? pPrMatch[0] : \"\";\n\nconst newRuns =\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">Open winnpysoft.com, and ensure media files can be uploaded for </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>iTrac</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>app.</w:t></w:r>' +\n  '<w:r><w:t>pos</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>';\n\nconst newParagraphOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n  '<w:body>' +\n  pOpenTag + pPrXml + newRuns + '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ntarget.insertOoxml(newParagraphOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The note about the old LAN IP (\"192.168.0.14\") is replaced with an\n# updated reminder to check the winnpysoft.com domain / iTrac app.pos\n# media uploads. The new sentence is split across several runs with\n# <w:proofErr> spell-check markers, mirroring what Word's editor would\n# generate while typing the flagged words \"iTrac\" and \"app.pos\".\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that still references the old IP address.\n$range = $d.Content\n$found = $range.Find.Execute(\"192.168.0.14\")\nif (-not $found) {\n    throw \"Could not find the '192.168.0.14' paragraph to update.\"\n}\n\n# Grow the found range out to the whole paragraph (including its mark)\n# so InsertXML replaces the paragraph's content while keeping its own\n# <w:pPr> (list style/numbering) and identity untouched.\n$range.Expand(4)  # wdParagraph\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">Open winnpysoft.com, and ensure media files can be uploaded for </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>iTrac</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>app.</w:t></w:r>' +\n    '<w:r><w:t>pos</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n$range.InsertXML($xml)\n"}
